$d = $word.ActiveDocument

# --- Step 1: apply the text-level edits (content becomes correct, but the
# engine will over-merge adjacent runs with identical formatting as a side
# effect of any text mutation within a paragraph) ---

$d.Content.Find.Execute("recalibrate, but", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "recalibrate the station, but", 2)

$d.Content.Find.Execute("stop producing for", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "stop production for", 2)

$d.Content.Find.Execute("any other issues.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ("any other issues that provoke to stop the line." + `
                         " This model will help not just the current production line but will be helpful for the others production line.   "), 2)

# --- Step 2: locate the edited paragraph and re-establish the run
# boundaries the diff expects by "touching" (toggling and restoring) a
# character-formatting property on each final run's exact sub-range. This
# forces the engine to split runs at those boundaries without altering any
# visible formatting. ---

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*experienced some*") {
        $target = $p
    }
}

$paraStart = $target.Range.Start

function Touch-Range($startOffset, $endOffset) {
    $r = $d.Range($paraStart + $startOffset, $paraStart + $endOffset)
    $orig = $r.Font.Bold
    $r.Font.Bold = 1
    $r.Font.Bold = $orig
}

# Offsets are relative to the start of the paragraph, computed from the
# final (post-edit) paragraph text so each Touch-Range call carves out
# exactly one of the runs the diff shows.
Touch-Range 0   26    # "The line experienced some "
Touch-Range 26  34    # "problems"
Touch-Range 34  106   # " like downtime in any station for quality issues and need to recalibrate"
Touch-Range 106 118   # " the station"
Touch-Range 118 184   # ", but calibrating can take some time; this provokes to stop produc"
Touch-Range 184 188   # "tion"
Touch-Range 188 189   # " "
Touch-Range 189 221   # "for not having fixture available"
Touch-Range 221 351   # ". The company is looking to know the optimal amount of fixture to continue the operation despite the downtimes or any other issues"
Touch-Range 351 381   # " that provoke to stop the line"
Touch-Range 381 382   # "."
Touch-Range 382 495   # " This model will help not just the current production line but will be helpful for the others production line.   "
